# ---------------------------------------------------------------------------
# Onderzoek-data.xlsx -- add a second address block (city2/street2/house_nr2/
# postalcode2) right after the existing address fields (city/street/house_nr/
# postalcode), pushing the remaining rows down, and update the sheet's
# view/selection state.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for 4 new rows: move the current contents of A15:A24
#    (created_at .. amount) down to A19:A28, working from the bottom up so
#    that we never overwrite a source cell before it has been read.
for ($r = 24; $r -ge 15; $r--) {
    $source = $ws.Cells.Item($r, 1)
    $target = $ws.Cells.Item($r + 4, 1)
    $target.Value = $source.Value2
}

# 2) Give the 4 freed-up rows (A15:A18) the same look as the other address
#    fields (city/street/house_nr/postalcode, rows 11-14) by copying their
#    formatting over before writing the new labels.
$ws.Cells.Item(14, 1).Copy()
$ws.Range("A15:A18").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# 3) Fill in the new field names.
$ws.Range("A15").Value = "city2"
$ws.Range("A16").Value = "street2"
$ws.Range("A17").Value = "house_nr2"
$ws.Range("A18").Value = "postalcode2"

# 4) Update the sheet view: scroll so row 7 is at the top and select C15.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
[void]$ws.Range("C15").Select()
